# Generate Report for Handoff
# Updates the localization-status report: the file
# "5083437c-5b8e-4167-95cd-74143c289871.md" (and the other files that were
# "Ready for handoff" with priority "ht") get a fresh handoff pass recorded:
#   - Overview!G        -> new "Latest HO Xliff Generate Date" timestamp
#   - zh-cn!E (Priority) ht -> mt, zh-cn!H (Latest Handoff Datetime) -> new timestamp
#   - de-de!E (Priority) ht -> mt, de-de!H (Latest Handoff Datetime) -> new timestamp

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13, 14, 16)

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-06 10:28:51"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "mt"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-06 10:28:45"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "mt"
    $dede.Cells.Item($r, 8).Value = "2016-09-06 10:28:51"
}
